$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D60").Value2 = "Work on homework for ninth lecture"
$ws.Range("D61").Formula = "=D60"
$v = $ws.Range("D61").Value2
Write-Host "D61 value:" $v
$f = $ws.Range("D61").Formula
Write-Host "D61 formula:" $f
